$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = '1,00'
$ws.Range("F3").Value = '4,00'
$ws.Range("F4").Value = '5,00'
$ws.Range("C5").Value = 'BAV24G0I1C'
$ws.Range("D5").Value = 'BALL VALVE W/INTEGRAL WELDED 2 NIPPLES, FB, FLOATING BALL, API 608, API 598, A105, CL 800, SW W/2 PE NIPPLES, MNF STD, SS316 BALL, SS316 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, API 607, LO'
$ws.Range("E5").Value = '0,75'
$ws.Range("F5").Value = '1,00'
$ws.Range("H5").Value = 'CSO'
$ws.Range("E6").Value = '8,00'
$ws.Range("F6").Value = '1,00'
$ws.Range("E7").Value = '10,00'
$ws.Range("F7").Value = '2,00'
$ws.Range("E8").Value = '14,00'
$ws.Range("F8").Value = '3,00'
$ws.Range("C9").Value = 'CKV21A0B2B'
$ws.Range("D9").Value = 'SWING CHECK VALVE FL, API 594, API 598, A216 GR WCB, CL 150, INST HORIZ/VERT, RF, B16.5, BOLTED COVER, SPW SS304/GRAPH, RENEWABLE SEATS, TRIM #8'
$ws.Range("E9").Value = '24,00'
$ws.Range("E10").Value = '3,00'
$ws.Range("F10").Value = '1,00'
$ws.Range("E11").Value = '6,00'
$ws.Range("F11").Value = '1,00'
$ws.Range("D12").Value = 'GATE VALVE FL, API 600, API 598, A216 GR WCB, CL 150, RF, B16.5, BB, SPW SS304/GRAPH, PKG GRAPH, TRIM #8, RENEWABLE SEATS, FLEXIBLE WEDGE, STEM OS&Y/RSNRO, HO'
$ws.Range("E12").Value = '8,00'
$ws.Range("F12").Value = '1,00'
$ws.Range("D13").Value = 'GATE VALVE FL, API 600, API 598, A216 GR WCB, CL 150, RF, B16.5, BB, SPW SS304/GRAPH, PKG GRAPH, TRIM #8, RENEWABLE SEATS, FLEXIBLE WEDGE, STEM OS&Y/RSNRO, HO'
$ws.Range("E13").Value = '8,00'
$ws.Range("F13").Value = '1,00'
$ws.Range("H13").Value = 'CSO'
$ws.Range("C14").Value = 'GAV21A0B2B'
$ws.Range("D14").Value = 'GATE VALVE FL, API 600, API 598, A216 GR WCB, CL 150, RF, B16.5, BB, SPW SS304/GRAPH, PKG GRAPH, TRIM #8, RENEWABLE SEATS, FLEXIBLE WEDGE, STEM OS&Y/RSNRO, HO'
$ws.Range("E14").Value = '12,00'
$ws.Range("F14").Value = '3,00'
$ws.Range("C15").Value = 'GAV21A0B2B'
$ws.Range("D15").Value = 'GATE VALVE FL, API 600, API 598, A216 GR WCB, CL 150, RF, B16.5, BB, SPW SS304/GRAPH, PKG GRAPH, TRIM #8, RENEWABLE SEATS, FLEXIBLE WEDGE, STEM OS&Y/RSNRO, GO'
$ws.Range("E15").Value = '14,00'
$ws.Range("F15").Value = '6,00'
$ws.Range("C16").Value = 'GAV21A0B2B'
$ws.Range("D16").Value = 'GATE VALVE FL, API 600, API 598, A216 GR WCB, CL 150, RF, B16.5, BB, SPW SS304/GRAPH, PKG GRAPH, TRIM #8, RENEWABLE SEATS, FLEXIBLE WEDGE, STEM OS&Y/RSNRO, GO'
$ws.Range("E16").Value = '16,00'
$ws.Range("F16").Value = '3,00'
$ws.Range("E17").Value = '2,00'
$ws.Range("F17").Value = '1,00'
$ws.Range("C18").Value = 'GAV413C3J2G'
$ws.Range("D18").Value = 'GATE VALVE FL, MSS SP-128 TYPE II, A536 Gr 65-45-12, CL 125, FF, B16.1, BB, NON METALLIC FLAT GASKET EPDM, PKG EPDM, SOLID WEDGE, WEDGE DUCTILE IRON EPDM ENCAPSULATED, STEM BRONZE, OS&Y/RSNRO, HO, UL LISTED/FM APPROVED'
$ws.Range("E18").Value = '3,00'
$ws.Range("F18").Value = '4,00'
$ws.Range("C19").Value = 'GAV413C3J2G'
$ws.Range("D19").Value = 'GATE VALVE FL, MSS SP-128 TYPE II, A536 Gr 65-45-12, CL 125, FF, B16.1, BB, NON METALLIC FLAT GASKET EPDM, PKG EPDM, SOLID WEDGE, WEDGE DUCTILE IRON EPDM ENCAPSULATED, STEM BRONZE, OS&Y/RSNRO, HO, UL LISTED/FM APPROVED'
$ws.Range("E19").Value = '6,00'
$ws.Range("F19").Value = '4,00'
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = '-'
$ws.Range("C20").Value = 'GAV413C3J2G'
$ws.Range("D20").Value = 'GATE VALVE FL, MSS SP-128 TYPE II, A536 Gr 65-45-12, CL 125, FF, B16.1, BB, NON METALLIC FLAT GASKET EPDM, PKG EPDM, SOLID WEDGE, WEDGE DUCTILE IRON EPDM ENCAPSULATED, STEM BRONZE, OS&Y/RSNRO, HO, UL LISTED/FM APPROVED'
$ws.Range("E20").Value = '8,00'
$ws.Range("F20").Value = '1,00'
$ws.Range("G20").Value = 'e.a'
$ws.Range("H20").Value = '-'
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = '-'
$ws.Range("C21").Value = 'GAV413C3J2G'
$ws.Range("D21").Value = 'GATE VALVE FL, MSS SP-128 TYPE II, A536 Gr 65-45-12, CL 125, FF, B16.1, BB, NON METALLIC FLAT GASKET EPDM, PKG EPDM, SOLID WEDGE, WEDGE DUCTILE IRON EPDM ENCAPSULATED, STEM BRONZE, OS&Y/RSNRO, HO, UL LISTED/FM APPROVED'
$ws.Range("E21").Value = '14,00'
$ws.Range("F21").Value = '6,00'
$ws.Range("G21").Value = 'e.a'
$ws.Range("H21").Value = '-'
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = '-'
$ws.Range("C22").Value = 'GAV413C3J2G'
$ws.Range("D22").Value = 'GATE VALVE FL, MSS SP-128 TYPE II, A536 Gr 65-45-12, CL 125, FF, B16.1, BB, NON METALLIC FLAT GASKET EPDM, PKG EPDM, SOLID WEDGE, WEDGE DUCTILE IRON EPDM ENCAPSULATED, STEM BRONZE, OS&Y/RSNRO, HO, UL LISTED/FM APPROVED'
$ws.Range("E22").Value = '16,00'
$ws.Range("F22").Value = '4,00'
$ws.Range("G22").Value = 'e.a'
$ws.Range("H22").Value = '-'
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = '-'
$ws.Range("C23").Value = 'GAV415J4J2G'
$ws.Range("D23").Value = 'HOSE GATE VALVE THD, MSS SP-80, B62 UNS C83600, 300 PSI CWP, FNPTxNH W/CAP AND CHAIN, FNPT B1.20.1 AND NH NFPA 1963, SCREWED BONNET, PKG NON ASBESTOS, SOLID WEDGE, WEDGE B62, SEAT&STEM BRONZE, NON RISING STEM, HO, UL LISTED/FM APPROVED'
$ws.Range("E23").Value = '2,5'
$ws.Range("F23").Value = '3,00'
$ws.Range("G23").Value = 'e.a'
$ws.Range("H23").Value = '-'
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = '-'
$ws.Range("C24").Value = 'GLV415J4J2G'
$ws.Range("D24").Value = 'ANGLE HOSE VALVE THD, MSS SP-80, B62 UNS C83600, 300 PSI CWP, FNPTxNH W/CAP AND CHAIN, FNPT B1.20.1 AND NH NFPA 1963, SCREW-IN BONNET, RENEWABLE DISC, PKG NON ASBESTOS, DISC&STEM BRONZE, RISING STEM, HO, UL LISTED/FM APPROVED'
$ws.Range("E24").Value = '2,5'
$ws.Range("F24").Value = '24,00'
$ws.Range("G24").Value = 'e.a'
$ws.Range("H24").Value = '-'
